$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3968

# Row 116
$ws.Range("H116").Value = 3272.4546
$ws.Range("I116").Value = 2351.8333
$ws.Range("K116").Value = 2351.8333
$ws.Range("M116").Value = 1090.1667

# Row 141
$ws.Range("H141").Value = 1396.6666
$ws.Range("I141").Value = 1076
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 3228
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 1952
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 58824572
$ws.Range("I61").Value = 90910030
$ws.Range("J61").Value = 1234.3334
$ws.Range("K61").Value = 90910030
$ws.Range("L61").Value = 1234.3334
$ws.Range("M61").Value = -90909818
$ws.Range("N61").Value = -1658.3334

# Row 122
$ws.Range("H122").Value = 1169.1034
$ws.Range("I122").Value = 794.7273
$ws.Range("K122").Value = 2384.1819
$ws.Range("M122").Value = 65.81809999999996

# Row 132
$ws.Range("H132").Value = 2696.516
$ws.Range("I132").Value = 2236.0833
$ws.Range("J132").Value = 4275.143
$ws.Range("K132").Value = 6708.249899999999
$ws.Range("L132").Value = 12825.429
$ws.Range("M132").Value = -4178.249899999999
$ws.Range("N132").Value = -17885.429

# Row 136
$ws.Range("H136").Value = 58824572
$ws.Range("I136").Value = 90910030
$ws.Range("J136").Value = 1234.3334
$ws.Range("K136").Value = 272730090
$ws.Range("L136").Value = 3703.0002
$ws.Range("M136").Value = -272727540
$ws.Range("N136").Value = -8803.0002

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4301.6
$ws.Range("I20").Value = 4254
$ws.Range("J20").Value = 4333.3335
$ws.Range("K20").Value = 4254
$ws.Range("L20").Value = 4333.3335
$ws.Range("M20").Value = -4007
$ws.Range("N20").Value = -4827.3335

# Row 22
$ws.Range("H22").Value = 1001
$ws.Range("I22").Value = 1001
$ws.Range("K22").Value = 1001
$ws.Range("M22").Value = -828

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 6898829
$ws.Range("I62").Value = 2352.1738
$ws.Range("J62").Value = 33335324
$ws.Range("K62").Value = 2352.1738
$ws.Range("L62").Value = 33335324
$ws.Range("M62").Value = -1728.1738
$ws.Range("N62").Value = -33336572

# Row 65
$ws.Range("H65").Value = 6898829
$ws.Range("I65").Value = 2352.1738
$ws.Range("J65").Value = 33335324
$ws.Range("K65").Value = 11760.869
$ws.Range("L65").Value = 166676620
$ws.Range("M65").Value = -8640.869000000001
$ws.Range("N65").Value = -166682860

# Row 99
$ws.Range("H99").Value = 1650
$ws.Range("I99").Value = 1627.6
$ws.Range("J99").Value = 1706
$ws.Range("K99").Value = 1627.6
$ws.Range("L99").Value = 1706
$ws.Range("M99").Value = -129.5999999999999
$ws.Range("N99").Value = -4702

# Row 126
$ws.Range("H126").Value = 1650
$ws.Range("I126").Value = 1627.6
$ws.Range("J126").Value = 1706
$ws.Range("K126").Value = 4882.799999999999
$ws.Range("L126").Value = 5118
$ws.Range("M126").Value = -2412.799999999999
$ws.Range("N126").Value = -10058

# Row 134
$ws.Range("H134").Value = 29414106
$ws.Range("I134").Value = 2555.5
$ws.Range("J134").Value = 166668000
$ws.Range("K134").Value = 7666.5
$ws.Range("L134").Value = 500004000
$ws.Range("M134").Value = -5131.5
$ws.Range("N134").Value = -500009070

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1294.1904
$ws.Range("I5").Value = 1411
$ws.Range("K5").Value = 4233
$ws.Range("M5").Value = -4121

# Row 131
$ws.Range("H131").Value = 18185136
$ws.Range("I131").Value = 90909450
$ws.Range("J131").Value = 4057.0908
$ws.Range("K131").Value = 272728350
$ws.Range("L131").Value = 12171.2724
$ws.Range("M131").Value = -272723310
$ws.Range("N131").Value = -22251.2724

# Row 132
$ws.Range("H132").Value = 840.25
$ws.Range("I132").Value = 795.3333
$ws.Range("J132").Value = 975
$ws.Range("K132").Value = 7157.9997
$ws.Range("L132").Value = 8775
$ws.Range("M132").Value = -4627.9997
$ws.Range("N132").Value = -13835

# Row 135
$ws.Range("H135").Value = 1294.1904
$ws.Range("I135").Value = 1411
$ws.Range("K135").Value = 12699
$ws.Range("M135").Value = -10164

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 37502970
$ws.Range("I70").Value = 41669770
$ws.Range("K70").Value = 41669770
$ws.Range("M70").Value = -41669500

# Row 73
$ws.Range("H73").Value = 37502970
$ws.Range("I73").Value = 41669770
$ws.Range("K73").Value = 41669770
$ws.Range("M73").Value = -41668834

# Row 80
$ws.Range("H80").Value = 6813.25
$ws.Range("I80").Value = 7726.25
$ws.Range("K80").Value = 7726.25
$ws.Range("M80").Value = -6728.25

# Row 83
$ws.Range("H83").Value = 6813.25
$ws.Range("I83").Value = 7726.25
$ws.Range("K83").Value = 38631.25
$ws.Range("M83").Value = -33639.25

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 837.1111
$ws.Range("I22").Value = 498
$ws.Range("J22").Value = 879.5
$ws.Range("K22").Value = 498
$ws.Range("L22").Value = 879.5
$ws.Range("M22").Value = -203
$ws.Range("N22").Value = -1469.5

# Row 27
$ws.Range("H27").Value = 837.1111
$ws.Range("I27").Value = 498
$ws.Range("J27").Value = 879.5
$ws.Range("K27").Value = 498
$ws.Range("L27").Value = 879.5
$ws.Range("M27").Value = -391
$ws.Range("N27").Value = -1093.5

# Row 46
$ws.Range("H46").Value = 4001.8462
$ws.Range("I46").Value = 763.3333
$ws.Range("J46").Value = 4973.4
$ws.Range("K46").Value = 763.3333
$ws.Range("L46").Value = 4973.4
$ws.Range("M46").Value = -575.3333
$ws.Range("N46").Value = -5349.4

# Row 82
$ws.Range("H82").Value = 2127.077
$ws.Range("I82").Value = 1968.3636
$ws.Range("K82").Value = 1968.3636
$ws.Range("M82").Value = -1607.3636

# Row 85
$ws.Range("H85").Value = 2127.077
$ws.Range("I85").Value = 1968.3636
$ws.Range("K85").Value = 1968.3636
$ws.Range("M85").Value = -720.3635999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 638.36365
$ws.Range("I107").Value = 553.6667
$ws.Range("J107").Value = 740
$ws.Range("K107").Value = 1661.0001
$ws.Range("L107").Value = 2220
$ws.Range("M107").Value = 258.9999
$ws.Range("N107").Value = -6060

# Row 136
$ws.Range("H136").Value = 914.075
$ws.Range("I136").Value = 859.13336
$ws.Range("J136").Value = 1078.9
$ws.Range("K136").Value = 2577.40008
$ws.Range("L136").Value = 3236.7
$ws.Range("M136").Value = -27.40008000000034
$ws.Range("N136").Value = -8336.700000000001
